# Inventario.xlsx edit script
# Summary of changes (per commit "Se agrego el documento con los integrntes del grupo XD"):
#  - productos: remove "Aguacate" row, bump "Waifus" existencia 5 -> 20,
#               add new product "Te Frio" (precio 4, existencia 10)
#  - cliente:   change Manolo Sandoval's nit 123456789 -> 789654123 and move
#               him to the bottom of the list, add new client "Marco Valdez"
#               (nit c/f, direccion Ciudad)
#  - pedido:    add new order row for Marco Valdez / Aguacate, cantidad 2, valor 8

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "productos"
# ---------------------------------------------------------------------------
$wsProductos = $wb.Worksheets.Item("productos")

# Remove the "Aguacate" row (row 2) - everything below shifts up one row.
$wsProductos.Rows(2).Delete()

# "Waifus" is now on row 11; update its existencia from 5 to 20.
$wsProductos.Range("C11").Value = 20

# Add the new product "Te Frio" on the next free row (row 12).
$wsProductos.Range("A12").Value = "Te Frio"
$wsProductos.Range("B12").Value = 4
$wsProductos.Range("C12").Value = 10

# ---------------------------------------------------------------------------
# Sheet "cliente"
# ---------------------------------------------------------------------------
$wsCliente = $wb.Worksheets.Item("cliente")

# Remove Manolo Sandoval's original row (row 2) - rest shift up one row.
$wsCliente.Rows(2).Delete()

# Re-add Manolo Sandoval at the bottom (row 7) with his updated nit, same city.
# Force the nit cell to Text format first so the numeric-looking string isn't
# reinterpreted as a number (matches the other "nit" cells in the column),
# then restore the Normal style so the cell doesn't carry leftover formatting.
$wsCliente.Range("A7").Value = "Manolo Sandoval"
$wsCliente.Range("B7").NumberFormat = "@"
$wsCliente.Range("B7").Value = "789654123"
$wsCliente.Range("B7").Style = "Normal"
$wsCliente.Range("C7").Value = "Aldea El Juez"

# Add the new client "Marco Valdez" on the next free row (row 8).
$wsCliente.Range("A8").Value = "Marco Valdez"
$wsCliente.Range("B8").Value = "c/f"
$wsCliente.Range("C8").Value = "Ciudad"

# ---------------------------------------------------------------------------
# Sheet "pedido"
# ---------------------------------------------------------------------------
$wsPedido = $wb.Worksheets.Item("pedido")

# Add the new order on the next free row (row 8).
$wsPedido.Range("A8").Value = "Marco Valdez"
$wsPedido.Range("B8").Value = "Aguacate"
$wsPedido.Range("C8").Value = 2
$wsPedido.Range("D8").Value = 8
